$d = $word.ActiveDocument

# Paragraph 2: "First Draft" -> apply strikethrough to the run and the paragraph mark
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.StrikeThrough = $true

# Paragraph 4: the last paragraph (with the _GoBack bookmark) -> insert "Second Draft" text
# before the bookmark, at the start of the paragraph.
$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range
$r.Collapse(1)
$r.InsertBefore("Second Draft")
